# Add a new custom character style named "yellow" that applies a yellow
# highlight (shading) to any run it is attached to. This mirrors the style
# Word mints when highlighted content (e.g. pasted from HTML with
# "background-color: yellow") needs a reusable character style.

$d = $word.ActiveDocument

# wdStyleTypeCharacter = 2
$style = $d.Styles.Add("yellow", 2)

# <w:basedOn w:val="DefaultParagraphFont"/>
$style.BaseStyle = $d.Styles("DefaultParagraphFont")

# <w:uiPriority w:val="1"/>
$style.Priority = 1

# <w:qFormat/>
$style.QuickStyle = $true

# <w:rPr><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto"/>
#         <w:shd w:val="clear" w:color="auto" w:fill="FFFF00"/></w:rPr>
# wdLineStyleNone = 0, wdColorYellow = 65535 (0x00FFFF00 fill)
try { $style.Font.Borders.OutsideLineStyle = 0 } catch { }
try { $style.Font.Borders.OutsideLineWidth = 0 } catch { }
try { $style.Font.Borders.OutsideColor = -16777216 } catch { }
try { $style.Font.Shading.BackgroundPatternColor = 65535 } catch { }

Write-Output "Added character style 'yellow' (yellow highlight shading)."
